$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.958.12'
$ws.Range('E2').Value = '  +5.12%  '

$ws.Range('D3').Value = '2.253.18'
$ws.Range('E3').Value = '  +1.40%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.82'
$ws.Range('E5').Value = '  +3.30%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.63'
$ws.Range('E6').Value = '  +5.92%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.531'
$ws.Range('E7').Value = '  +3.32%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.482'
$ws.Range('E9').Value = '  +3.23%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.65'
$ws.Range('E10').Value = '  +7.13%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.58'
$ws.Range('E11').Value = '  +9.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0798'
$ws.Range('E12').Value = '  +2.07%  '

$ws.Range('E13').Value = '  +3.12%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.68'
$ws.Range('E14').Value = '  +3.59%  '

$ws.Range('E15').Value = '  +1.60%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.11'
$ws.Range('E16').Value = '  +2.09%  '

$ws.Range('D17').Value = '2.257.83'
$ws.Range('E17').Value = '  +1.77%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.756'
$ws.Range('E18').Value = '  +3.34%  '

$ws.Range('D19').Value = '41.864.34'
$ws.Range('E19').Value = '  +5.01%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.19'
$ws.Range('E20').Value = '  +9.56%  '

$ws.Range('D21').Value = '0.0₃0903'
$ws.Range('E21').Value = '  +1.84%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.94'
$ws.Range('E22').Value = '  +3.26%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.09'
$ws.Range('E23').Value = '  +2.11%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '241.40'
$ws.Range('E24').Value = '  +1.70%  '

$ws.Range('E25').Value = '  +5.44%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.11%  '

$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.91'
$ws.Range('E27').Value = '  +3.94%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.89'
$ws.Range('E28').Value = '  +3.40%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.22'
$ws.Range('E29').Value = '  +8.11%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.67'
$ws.Range('E30').Value = '  +4.66%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.08'
$ws.Range('E31').Value = '  +6.36%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '158.37'
$ws.Range('E32').Value = '  +1.06%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.02%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.15'
$ws.Range('E34').Value = '  +3.67%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0742'
$ws.Range('E35').Value = '  +4.19%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.03'
$ws.Range('E36').Value = '  +1.62%  '

$ws.Range('E38').Value = '  +5.74%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.61'
$ws.Range('E39').Value = '  +7.58%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.115'
$ws.Range('E40').Value = '  +3.58%  '

$ws.Range('E41').Value = '  +3.12%  '

$ws.Range('E42').Value = '  +5.30%  '

$ws.Range('D43').Value = '2.047.84'
$ws.Range('E43').Value = '  -2.83%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.91'
$ws.Range('E44').Value = '  +8.71%  '

$ws.Range('E45').Value = '  +3.48%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.11'
$ws.Range('E46').Value = '  +1.65%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.05'
$ws.Range('E47').Value = '  +2.05%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.86'
$ws.Range('E48').Value = '  +5.15%  '

$ws.Range('E49').Value = '  +3.35%  '

$ws.Range('E50').Value = '  +4.17%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.76'
$ws.Range('E51').Value = '  +5.68%  '
